# Requisitos.xlsx - "docu excel + arreglo índice"
#
# 1) Porcentaje completado (G7, G8) updated from 50% to 75%.
# 2) A new task row (row 19) is (re)populated describing the new
#    "Mostrar Contraseña" feature - a completed, low priority task.
# 3) Sheet view scroll position / selection refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Porcentaje completado updates -----------------------------------
$ws.Range("G7").Value = 0.75
$ws.Range("G8").Value = 0.75

# --- Row 19: "Mostrar Contraseña" task --------------------------------
$ws.Range("B19").Value = "Mostrar Contraseña"
$ws.Range("C19").Value = "Baja"
$ws.Range("D19").Value = "Completada"
$ws.Range("G19").Value = 1
$ws.Range("I19").Value = "Muestra la contraseña de todos los formularios."
$ws.Range("J19").Value = "OK"

# --- Sheet view: scroll / selection ------------------------------------
$ws.Range("G8").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
